$wb = $excel.ActiveWorkbook

# Work on the "config" sheet (last sheet in the workbook)
$ws = $wb.Worksheets.Item("config")

# Insert a new "commodity" column (C) with a header, shifting the old
# "level"/"unit" data. Final layout:
#   Row1: node | sector | commodity | level
#   Row2: R11_AFR | i_therm | i_therm | useful
#   Row3: R11_CPA | rc_spec | i_therm | useful

$ws.Range("C1").Value = "commodity"
$ws.Range("D1").Value = "level"

# Fit column C to its new header before filling in the shorter data rows,
# matching the width Excel would have settled on while the column only
# contained the header text.
[void]$ws.Columns.Item(3).AutoFit()

$ws.Range("C2").Value = "i_therm"
$ws.Range("D2").Value = "useful"

$ws.Range("C3").Value = "i_therm"
$ws.Range("D3").Value = "useful"

# Remove the now-unused row 4 (previously held the year value 2040)
[void]$ws.Range("A4:D4").Clear()

[void]$ws.Range("D4").Select()

# Adjust the selection on the MERtoPPP sheet as well
$ws2 = $wb.Worksheets.Item("MERtoPPP")
[void]$ws2.Range("D13").Select()
[void]$ws.Activate()
